# Applies the "Add files via upload" commit:
#  - Appends 5 new participant rows (297-301) to the 'Отобранные участники'
#    sheet, copying number formats from existing rows.
#  - Extends the shared I-column formula down through the new rows.
#  - Extends the sheet's AutoFilter range and the _FilterDatabase defined
#    name from $A$1:$AA$296 to $A$1:$AA$301.
#  - The whole-column COUNTIF summary cells (X4/X5/X7) recompute on their
#    own once the new "Да" values exist in columns J/K/M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New row data (A..U), values taken from the target XLSX. $null means
#    the source cell is blank.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=297; A="Чепелева Мария Антоновна";         B="@mrukly";      C=89513761868; D="Android"; E="Новосибирская область"; F="Новосибирск ";      G="Женский"; H=21;
       J="Да";    K="Запас"; L="Запас"; M="Да";    N="Сценарий 2"; O="Без обращения"; P=$null;        Q=$null;        R=$null; S=$null; T="Без обращения"; U="Сценарий 4" },
    @{ Row=298; A="Кретов Тимофей Николаевич";        B="@tkretov";     C=89652282950; D="Android"; E="Москва";                F="Москва";              G="Мужской"; H=39;
       J="Запас"; K="Запас"; L="Запас"; M="Да";    N=$null;        O=$null;           P=$null;        Q=$null;        R=$null; S=$null; T="Сценарий 2";     U="Без обращения" },
    @{ Row=299; A="Крамаренко Анна Ивановна";         B="@miasamr";     C=89514944126; D="iOS";     E="Ростовская об";         F="Ростов-на-Дону";      G="Женский"; H=21;
       J="Запас"; K="Да";    L="Запас"; M="Да";    N=$null;        O=$null;           P="Сценарий 2"; Q="Сценарий 5"; R=$null; S=$null; T="Без обращения"; U="Сценарий 4" },
    @{ Row=300; A="Макарычева Светлана Геннадьевна";  B="@SvetlaniaM";  C=89035825849; D="Android"; E="Россия";                F="Москва";              G="Женский"; H=45;
       J="Да";    K="Нет, не являюсь клиентом данного банка"; L="Нет, не являюсь клиентом данного банка"; M="Нет, не являюсь клиентом данного банка";
       N="Сценарий 2"; O="Без обращения"; P=$null; Q=$null; R=$null; S=$null; T=$null; U=$null },
    @{ Row=301; A="Пошина Мария Викторовна";          B="@Ruandil";     C=89241177290; D="Android"; E="Хабаровский край ";     F="Хабаровск ";          G="Женский"; H=50;
       J="Запас"; K="Нет, не являюсь клиентом данного банка"; L="Нет, не являюсь клиентом данного банка"; M="Да";
       N=$null; O=$null; P=$null; Q=$null; R=$null; S=$null; T="Без обращения"; U="Сценарий 4" }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

# ---------------------------------------------------------------------
# 2. Copy the row-296 number formats down for rows 297-300, and the
#    row-295 formats (which already carries the highlighted "I" style
#    used by row 301) down for row 301, then fill in the values/formula.
# ---------------------------------------------------------------------
$ws.Range("A296:U296").Copy() | Out-Null
$ws.Range("A297:U300").PasteSpecial(-4122) | Out-Null

$ws.Range("A295:U295").Copy() | Out-Null
$ws.Range("A301:U301").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

foreach ($r in $newRows) {
    $rowNum = $r.Row

    foreach ($col in $cols) {
        if ($col -eq "I") { continue }
        $val = $r[$col]
        $cell = $ws.Range($col + $rowNum)
        if ($val -eq $null) {
            $cell.Value = ""
        } else {
            $cell.Value = $val
        }
    }

    # Column I: same shared formula pattern as I290:I296, extended down.
    $ws.Cells.Item($rowNum, 9).Formula = '=COUNTIF(J' + $rowNum + ':M' + $rowNum + ', "Да, являюсь клиентом банка более полугода") + COUNTIF(J' + $rowNum + ':M' + $rowNum + ', "Да, являюсь клиентом банка менее полугода")'
}

# ---------------------------------------------------------------------
# 3. Extend the AutoFilter range from $A$1:$AA$296 to $A$1:$AA$301.
#    Toggling off first avoids Range.AutoFilter() acting as a pure
#    on/off switch when a filter is already active.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AA301").AutoFilter() | Out-Null

# ---------------------------------------------------------------------
# 4. Update the _xlnm._FilterDatabase defined name to match.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='" + $ws.Name + "'!`$A`$1:`$AA`$301"
    }
}

# ---------------------------------------------------------------------
# 5. Recalculate so the whole-column COUNTIF summaries (X4, X5, X7) and
#    the new I-column formulas carry fresh cached values.
# ---------------------------------------------------------------------
$excel.CalculateFull()

Write-Output "done"
